$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.328.37"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "2.031.70"
$ws.Range("E3").Value = "  +3.92%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.88"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0810"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.23"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.23%  "
$ws.Range("D13").Value = "2.326.81"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.851"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "2.025.96"
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "37.267.32"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +5.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.39"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("E28").Value = "  -5.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.99"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.04%  "
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0675"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +10.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.56"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.63"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.61%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0973"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.75"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.58"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.21%  "
$ws.Range("D47").Value = "1.379.87"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("E49").Value = "  +15.80%  "
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.34%  "
